# Applies the "Celestial Choreography" -> "History" rewrite described in the
# commit diff ("Three sets of 1000 word documents with three font types").
#
# Strategy: literal (non-wildcard) Find/Replace scoped to $d.Content for
# each sentence/fragment, in document order, so repeated punctuation (the
# lone "." runs between sentences) never causes an ambiguous match. Where
# the diff adds brand-new sentences/runs that did not exist before, the
# freshly rewritten Range is collapsed to its end and Range.InsertAfter is
# used to graft the new text on.

$d = $word.ActiveDocument

function Replace-Literal([string]$old, [string]$new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Replace-Literal: could not find [$old]"
    }
    return $rng
}

# ---------------------------------------------------------------------
# Title
# ---------------------------------------------------------------------
Replace-Literal "Celestial Choreography: A Journey Through Cosmic Motions" "History: A Tapestry of Human Endeavors" | Out-Null

# ---------------------------------------------------------------------
# Author name: "Anastasia Leontiev" -> "Dr" / "." / " Henry Longfellow"
# ---------------------------------------------------------------------
$rng = Replace-Literal "Anastasia Leontiev" "Dr"
$rng.Collapse(0)
$rng.InsertAfter(".")
$rng.Collapse(0)
$rng.InsertAfter(" Henry Longfellow")

# ---------------------------------------------------------------------
# Email line: "astrid" . "leontiev@earthnet" . "net" ->
#             "hlongfellow@schoolhistory" . "org"
# ---------------------------------------------------------------------
Replace-Literal "astrid" "hlongfellow@schoolhistory" | Out-Null
Replace-Literal "leontiev@earthnet" "org" | Out-Null
Replace-Literal ".org.net" ".org" | Out-Null

# ---------------------------------------------------------------------
# Body paragraph (sz 24), sentence by sentence
# ---------------------------------------------------------------------
Replace-Literal `
    "In the vast expanse of the universe, amidst the myriad celestial bodies, an intricate dance unfolds, revealing the harmony and elegance of cosmic motions" `
    "History has always enticed us, capturing our imagination with its tales of triumphs, tragedies, remarkable achievements, and lessons learned" | Out-Null

Replace-Literal `
    " From the grand ballet of planets orbiting stars to the  gravitational waltz of galaxies, the cosmos is filled with captivating movements that have captivated scientists, philosophers, and artists alike" `
    " Within its vast expanse, we unearth the foundations of civilizations, the intricacies of diverse cultures, and the evolution of human thought and action" | Out-Null

$rng = Replace-Literal `
    " Understanding these celestial choreographies provides invaluable insights into the fundamental laws of physics, the origins of our universe, and the captivating beauty of our place within it" `
    " It serves as a mosaic composed of countless individual stories, each contributing its unique hue to the broader canvas of our shared heritage"
$rng.Collapse(0)
$rng.InsertAfter(".")
$rng.Collapse(0)
$rng.InsertAfter(" As we delve into the depths of historical studies, we assume the role of explorers, unraveling the intricate threads that bind us to past events and shedding light on the origins of our present world")

Replace-Literal `
    "Each celestial body, whether a planet, star, or galaxy, possesses an inherent angular momentum, a conserved quantity that governs its rotation and revolution around a central axis or an external gravitational center" `
    "The tapestry of history is woven with the endeavors of countless individuals, both renowned and forgotten" | Out-Null

Replace-Literal `
    " This angular momentum plays a crucial role in determining the trajectories, shapes, and dynamics of celestial objects, weaving them into intricate patterns of perpetual motion" `
    " Kings and queens, generals and statesmen, revolutionaries and reformers, artists and thinkers - each has left their indelible mark on the canvas of human progress" | Out-Null

$rng = Replace-Literal `
    " These patterns, governed by the laws of gravitation and conservation of energy, exhibit a level of predictability and order, hinting at the underlying mathematical harmony of the universe" `
    " But history would be incomplete if it solely focused on the exploits of the few; it is the collective actions and contributions of ordinary individuals, the unsung heroes of history, that truly shape its course"
$rng.Collapse(0)
$rng.InsertAfter(".")
$rng.Collapse(0)
$rng.InsertAfter(" As students of history, we have the privilege of delving into the lives of these individuals, understanding their motivations, struggles, and triumphs, and appreciating the profound impact they have had on our world")

Replace-Literal `
    "The gravitational forces exerted between celestial bodies, in accordance with Newton's laws of motion, choreograph their graceful interactions" `
    "Furthermore, history offers us a lens through which we can examine ourselves, our values, and the challenges we confront in our ever-changing society" | Out-Null

Replace-Literal `
    " The interplay of attraction and repulsion, governed by the delicate balance between mass, velocity, and gravitational constant, orchestrates the celestial dance" `
    " By understanding the past, we gain perspective on the present and can make more informed decisions about the future" | Out-Null

$rng = Replace-Literal `
    " Planets take up elliptical paths around stars, moons encircle planets, and galaxies spiral in elegant synchrony, all guided by the invisible yet powerful forces that shape the fabric of spacetime" `
    " History compels us to confront uncomfortable truths, to grapple with difficult lessons, and to cultivate empathy for those who have come before us"
$rng.Collapse(0)
$rng.InsertAfter(".")
$rng.Collapse(0)
$rng.InsertAfter(" In this way, history serves as a guide, helping us navigate the complexities of human existence and equipping us with the knowledge and wisdom to shape a better future for ourselves and generations to come")

# ---------------------------------------------------------------------
# Summary heading: unchanged
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# Summary paragraph
# ---------------------------------------------------------------------
Replace-Literal `
    "The ballet of celestial motions is a mesmerizing spectacle that reveals the profound elegance and mathematical harmony of the universe" `
    "History, in its vastness, is a narrative of human experiences, triumphs and challenges, that weaves together the fabric of our past to illuminate the present and guide the future" | Out-Null

$rng = Replace-Literal `
    " From the intricate dance of planets orbiting stars to the gravitational waltz of galaxies, each celestial choreography holds secrets about the fundamental laws of physics, the origins of our universe, and the captivating interplay of forces that shape cosmic structures" `
    " Through "
$rng.Collapse(0)
$rng.InsertAfter("history, we gain insights into the lives of individuals who have shaped our world, the lessons learned from their endeavors, and the complex dynamics that have influenced the course of civilization")

Replace-Literal `
    " Understanding these cosmic movements enhances our " `
    " History is not merely a collection of facts; it is a tapestry of stories that holds the keys to our collective identity, values, and aspirations" | Out-Null

$rng = Replace-Literal "appreciation for the intricate beauty of the universe, inspiring awe and wonder in the face of the vast and mysterious cosmos" "."
$rng.Collapse(0)
$rng.InsertAfter(" By studying history, we embark on a transformative journey that deepens our understanding of ourselves and our place within the grand tapestry of human existence")

# ---------------------------------------------------------------------
# Trailing empty paragraph at the end of the body
# ---------------------------------------------------------------------
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()

Write-Host "done"
